$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append below the existing header + English row
$data = @(
    @("Catalan",   "50 (1-6), 100 (1-6), 200 (1-6), 300 (1-6), 600 (1-6)", "10.5281/zenodo.17352129"),
    @("Bosnian",   "50 (1-6), 100 (1-6), 200 (1-6), 300 (1-6), 600 (1-6)", "10.5281/zenodo.17344027"),
    @("Afrikaans", "50 (1-6), 100 (1-6), 200 (1-6), 300 (1-6), 600 (1-6)", "10.5281/zenodo.17328169"),
    @("Arabic",    "50 (1-6), 100 (1-6), 200 (1-6), 300 (1-4), 300_5_cbow", "10.5281/zenodo.17334562"),
    @("Arabic",    "300_5_sg, 300 (6), 500 (1-6)", "10.5281/zenodo.17334562"),
    @("French",    "50 (1-6), 100 (1-6), 200 (1-6), 300 (1-6), 600 (1-6)", "10.5281/zenodo.17337550"),
    @("Galician",  "50 (1-6), 100 (1-6), 200 (1-6), 300 (1-6), 600 (1-6)", "10.5281/zenodo.17343732")
)

$startRow = 3
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# Apply the built-in "Hyperlink" style to the French DOI cell (C8)
$ws.Range("C8").Style = "Hyperlink"

# Resize column B to fit the widest new content (matches Excel's computed
# "best fit" width of 44.5 for the longest file-list string in Aptos Narrow 12)
$ws.Columns.Item(2).ColumnWidth = 43.6

# Update the active selection to reflect the new extent, like Excel would after entry
$ws.Range("A10").Select()
